# Split the trailing "<number>}" / "<number>},"  token at the end of
# several JSON-snippet paragraphs into two runs: the numeric part (plus
# any trailing space) gets colored red, and the first one additionally
# gets a yellow highlight; the closing brace/comma stays in the
# original (unformatted) run style.

$d = $word.ActiveDocument

# Each entry: 1-based paragraph index, the literal text to search for
# within that paragraph, how many leading characters become the
# "red" run, and whether that red run should also be yellow-highlighted.
$targets = @(
    @{ Para = 7;  Text = "2003},"; SplitAt = 4; Highlight = $true  },
    @{ Para = 8;  Text = "1000},"; SplitAt = 4; Highlight = $false },
    @{ Para = 9;  Text = "2002 }"; SplitAt = 5; Highlight = $false },
    @{ Para = 19; Text = "1502 }"; SplitAt = 5; Highlight = $false },
    @{ Para = 25; Text = "502 }";  SplitAt = 4; Highlight = $false },
    @{ Para = 31; Text = "1022 }"; SplitAt = 5; Highlight = $false },
    @{ Para = 39; Text = "1503 }"; SplitAt = 5; Highlight = $false },
    @{ Para = 44; Text = "503 }";  SplitAt = 4; Highlight = $false },
    @{ Para = 48; Text = "533 }";  SplitAt = 4; Highlight = $false },
    @{ Para = 53; Text = "1033 }"; SplitAt = 5; Highlight = $false },
    @{ Para = 54; Text = "2033 }"; SplitAt = 5; Highlight = $false }
)

foreach ($t in $targets) {
    $paraRange = $d.Paragraphs.Item($t.Para).Range
    $found = $paraRange.Find.Execute($t.Text, $true, $false, $false, $false,
                                      $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: para $($t.Para) text '$($t.Text)'"
        continue
    }

    $matchStart = $paraRange.Start
    $matchEnd = $paraRange.End
    $splitPoint = $matchStart + $t.SplitAt

    $firstRange = $d.Range($matchStart, $splitPoint)
    $firstRange.Font.Color = 255
    if ($t.Highlight) {
        $firstRange.Font.HighlightColorIndex = 7
    }
}
